$d = $word.ActiveDocument

$d.Content.Find.Execute("76÷4=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=24, 1", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷8=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=38, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=21, 1", 2) | Out-Null
$d.Content.Find.Execute("30÷9=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("95÷3=31, 2", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=10, 1", 2) | Out-Null
$d.Content.Find.Execute("80÷6=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=17, 3", 2) | Out-Null
$d.Content.Find.Execute("58÷9=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2) | Out-Null
$d.Content.Find.Execute("21÷5=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2) | Out-Null
$d.Content.Find.Execute("12÷7=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=7, 2", 2) | Out-Null
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "25÷6=4, 1", 2) | Out-Null
$d.Content.Find.Execute("85÷7=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=5, 6", 2) | Out-Null
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "38÷9=4, 2", 2) | Out-Null
$d.Content.Find.Execute("55÷3=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷3=26, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷8=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2) | Out-Null
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "64÷8=8, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷8=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷6=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "97÷8=12, 1", 2) | Out-Null
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=10, 5", 2) | Out-Null
$d.Content.Find.Execute("64÷2=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=12, 4", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=11, 4", 2) | Out-Null
$d.Content.Find.Execute("78÷3=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=6, 1", 2) | Out-Null
$d.Content.Find.Execute("17÷6=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=19, 3", 2) | Out-Null
$d.Content.Find.Execute("46÷2=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷7=6, 1", 2) | Out-Null
$d.Content.Find.Execute("77÷5=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "87÷7=12, 3", 2) | Out-Null
